# SwaadSutra Daily Orders update - 2026-01-13T11:15:13.638Z
# New order placed by Ajay Dwarkunde (Pohe x1) - insert as the newest row
# on top of "Daily Orders", add matching line to "Items Breakdown", and
# refresh the "Summary" roll-up numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily Orders - insert a new row under the header for the new order
# ---------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Daily Orders")
$orders.Rows.Item(2).Insert()

# Inserting a row re-serialises the still-empty cells that shifted down
# (old rows 2 & 3, now 3 & 4) through the shared-string table, which
# leaves them holding an empty string instead of being truly blank.
# Re-clear them so they stay blank, matching the untouched source rows.
$orders.Range("E3").Value = ""
$orders.Range("L3").Value = ""
$orders.Range("M3").Value = ""
$orders.Range("N3").Value = ""
$orders.Range("E4").Value = ""
$orders.Range("L4").Value = ""
$orders.Range("M4").Value = ""
$orders.Range("N4").Value = ""

$orders.Range("A2").Value = 3
$orders.Range("B2").Value = "2026-01-13 11:15"
$orders.Range("C2").Value = "Ajay Dwarkunde"
$orders.Range("D2").Value = "b-703"

# Force text formatting on the phone number so Excel keeps it as the
# literal digit string instead of silently coercing it to a number.
$orders.Range("E2").NumberFormat = "@"
$orders.Range("E2").Value = "8087172173"

$orders.Range("F2").Value = "Pohe x1"
$orders.Range("G2").Value = 30
$orders.Range("H2").Value = "NEW"
$orders.Range("I2").Value = "PENDING"

# Force text formatting on the date-only cell so Excel keeps the literal
# "2026-01-13" string instead of silently coercing it to a date serial.
$orders.Range("J2").NumberFormat = "@"
$orders.Range("J2").Value = "2026-01-13"

$orders.Range("K2").Value = "18:50"
$orders.Range("L2").Value = ""
$orders.Range("M2").Value = ""
$orders.Range("N2").Value = ""

# ---------------------------------------------------------------------
# 2) Items Breakdown - insert a new row for Pohe above Wheat Chapati
# ---------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")
$items.Rows.Item(2).Insert()

$items.Range("A2").Value = "Pohe"
$items.Range("B2").Value = 1
$items.Range("C2").Value = 30

# ---------------------------------------------------------------------
# 3) Summary - refresh the totals to account for the new order
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A2").Value = 3
$summary.Range("B2").Value = 2
$summary.Range("G2").Value = 75
